$d = $word.ActiveDocument

# Remove the three runs of text around the "_GoBack" bookmark while
# keeping the bookmark itself, and remove the trailing sentence after it.
$d.Content.Find.Execute(" Als Userschnittstelle verwenden wir ein ganz einfache Terminal UI", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)
$d.Content.Find.Execute(" wo man anhand die Zeit, Film und Sitzplatz eine «Ticket» reservieren kann.", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)
